$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, shifting rows 116:203 down to 117:204
$ws.Rows.Item(116).Insert()

# Write the new row 116 values
$ws.Cells.Item(116, 1).Value = 8
$ws.Cells.Item(116, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(116, 3).Value = "Coquimbo"
$ws.Cells.Item(116, 4).Value = 44978
$ws.Cells.Item(116, 5).Value = 4
$ws.Cells.Item(116, 6).Value = 100112044
$ws.Cells.Item(116, 7).Value = "Perejil"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 2000
$ws.Cells.Item(116, 11).Value = 2300
$ws.Cells.Item(116, 12).Value = 2500
$ws.Cells.Item(116, 13).Value = 2400
$ws.Cells.Item(116, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(116, 16).Value = 1600
$ws.Cells.Item(116, 17).Value = 1.5
$ws.Cells.Item(116, 18).Value = "Hortaliza"
